$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(" Nov 2 2020", " Abu Dhabi", "Capitals won by 6 wickets (with 6 balls remaining)", "Royal Challengers Bangalore", "Delhi Capitals", "Washington Sundar ", "0", "1", "0", "0", "0.00"),
    @(" Nov 6 2020", " Abu Dhabi", "Sunrisers won by 6 wickets (with 2 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", "Washington Sundar ", "5", "6", "0", "0", "83.33"),
    @(" Oct 15 2020", " Sharjah", "Kings XI won by 8 wickets", "Royal Challengers Bangalore", "Kings XI Punjab", "Washington Sundar ", "13", "14", "1", "0", "92.85"),
    @(" Oct 25 2020", " Dubai (DSC)", "Super Kings won by 8 wickets (with 8 balls remaining)", "Royal Challengers Bangalore", "Chennai Super Kings", "Washington Sundar ", "5", "2", "1", "0", "250.00"),
    @(" Sep 24 2020", " Dubai (DSC)", "Kings XI won by 97 runs", "Royal Challengers Bangalore", "Kings XI Punjab", "Washington Sundar ", "30", "27", "2", "1", "111.11"),
    @(" Oct 31 2020", " Sharjah", "Sunrisers won by 5 wickets (with 35 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", "Washington Sundar ", "21", "18", "2", "0", "116.66"),
    @(" Oct 28 2020", " Abu Dhabi", "Mumbai won by 5 wickets (with 5 balls remaining)", "Royal Challengers Bangalore", "Mumbai Indians", "Washington Sundar ", "10", "6", "1", "0", "166.66"),
    @(" Oct 10 2020", " Dubai (DSC)", "RCB won by 37 runs", "Royal Challengers Bangalore", "Chennai Super Kings", "Washington Sundar ", "10", "10", "0", "1", "100.00"),
    @(" Oct 5 2020", " Dubai (DSC)", "Capitals won by 59 runs", "Royal Challengers Bangalore", "Delhi Capitals", "Washington Sundar ", "17", "11", "3", "0", "154.54")
)

# Force the whole data range to be stored as text, matching the
# numberStoredAsText / t="str" semantics of the original sheet.
$dataRange = $ws.Range("A2:K10")
$dataRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}

$ws.Range("A1:K10").NumberFormat = "@"
